$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# The "FilesTab" row's query (cell B4) is replaced with an updated Cypher
# query that joins through a parent node and formats the file size into a
# human readable unit (Bytes/KB/MB/GB/TB) instead of returning raw bytes.
$newQuery = "MATCH (f:file)-->(parent)" + [char]10 + `
"MATCH (f)-[:file_of_sample]->(samp)" + [char]10 + `
"MATCH (samp)-[:sample_of_study_subject]->(ss)" + [char]10 + `
"MATCH (ss)-[:study_subject_of_study]->(s)" + [char]10 + `
"MATCH (s)-[:study_of_program]->(p)" + [char]10 + `
"MATCH (d)-[:diagnosis_of_study_subject]->(ss)" + [char]10 + `
"MATCH (tp)-[:tp_of_diagnosis]->(d)" + [char]10 + `
"WHERE tp.chemotherapy_regimen IN [`"Dose dense AC (2 week cycles)`"]" + [char]10 + `
"WITH" + [char]10 + `
"        f, parent,p, ss, d,tp, s, samp," + [char]10 + `
"        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units," + [char]10 + `
"        toInteger(floor(log(f.file_size)/log(1024))) as i," + [char]10 + `
"        2 as precision" + [char]10 + `
"WITH" + [char]10 + `
"        f, parent,p, ss, d,tp, s, samp," + [char]10 + `
"        f.file_size /(1024^i) AS value," + [char]10 + `
"        10^precision AS factor," + [char]10 + `
"        units[i] as unit" + [char]10 + `
"WITH" + [char]10 + `
"        f, parent,p, ss, d,tp, s, samp, unit," + [char]10 + `
"        round(factor * value)/factor AS size" + [char]10 + `
"RETURN Distinct" + [char]10 + `
"    f.file_name AS ``File Name``," + [char]10 + `
"    head(labels(samp)) AS ``Association``," + [char]10 + `
"    f.file_description AS ``Description``," + [char]10 + `
"    f.file_format AS ``File Format``," + [char]10 + `
"     CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size," + [char]10 + `
"    p.program_acronym AS ``Program Code``," + [char]10 + `
"    s.study_acronym AS ``Arm``," + [char]10 + `
"    ss.study_subject_id AS ``Case ID``," + [char]10 + `
"    samp.sample_id AS ``Sample ID``" + [char]10 + `
"    order by f.file_name"

$ws.Range("B4").Value = $newQuery

# Match the author's workflow of clicking into the cell they just edited.
[void]$ws.Range("B4").Select()
